$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.85 = 6898.56 pesos`n✅ 6898.56 pesos = 1.85 = 885.64 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the "tasas" sheet numeric cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 539.997
$ws2.Range("O10").Value = 3725.2
$ws2.Range("N12").Value = 3738.94
$ws2.Range("O12").Value = 480.005
